# Add a "% of Q Drop's" column (I) to the Summer2013 GE master sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column.
$header = @{1 = "% of Q Drop's"}

# Data rows: default "0.00%" for every course/professor row, except
# OCNG-251 / WORMUTH J (row 43) which gets "3.57%".
$values = @{
    3  = "0.00%"
    6  = "0.00%"
    9  = "0.00%"
    12 = "0.00%"
    15 = "0.00%"
    18 = "0.00%"
    19 = "0.00%"
    22 = "0.00%"
    23 = "0.00%"
    26 = "0.00%"
    27 = "0.00%"
    30 = "0.00%"
    33 = "0.00%"
    36 = "0.00%"
    37 = "0.00%"
    40 = "0.00%"
    43 = "3.57%"
}

foreach ($r in $header.Keys) {
    $cell = $ws.Cells.Item([int]$r, 9)
    $cell.NumberFormat = "@"
    $cell.Value = $header[$r]
}

foreach ($r in $values.Keys) {
    $cell = $ws.Cells.Item([int]$r, 9)
    $cell.NumberFormat = "@"
    $cell.Value = $values[$r]
}

Write-Output "Added % of Q Drop's column (I) with $($values.Count) data rows"
